# BOM adjustments to R28, R29, config resistors
#
# Row 8 on the "bitaxeHex" sheet is the "10k" resistor BOM line. R28 and R29
# were mistakenly given their own line (row 22, value "10K", wrong DK/PARTNO
# pulled from the 5.6k line above). Fold R28, R29 into the correct "10k" row
# (row 8) and remove the now-redundant standalone row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bitaxeHex")

# 1. Update the "10k" resistor row (row 8): bump the qty and insert
#    "R28, R29" into the sorted Reference(s) list.
$ws.Range("B8").Value = 24
$ws.Range("C8").Value = "R3, R7, R15, R18, R22, R28, R29, R30, R47, R48, R59, R60, R67, R68, R75, R76, R81, R82, R86, R87, R100, R101, R105, R106"

# 2. Remove the old standalone "R28, R29 / 10K" row (row 22) entirely,
#    shifting everything below it up by one row.
$ws.Rows.Item(22).Select()
$ws.Rows.Item(22).Delete()
